$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-40 (A: key, B: localized value) to reflect new UI strings
# for dialog / options / confirm / hud elements, and reorder some rows.
$ws.Cells.Item(2,1).Value = "welcome"
$ws.Cells.Item(2,2).Value = "Welcome!"
$ws.Cells.Item(3,1).Value = "title"
$ws.Cells.Item(3,2).Value = "PERFECT CELL"
$ws.Cells.Item(4,1).Value = "none"
$ws.Cells.Item(4,2).Value = "None"
$ws.Cells.Item(5,1).Value = "test1"
$ws.Cells.Item(5,2).Value = "Test 1"
$ws.Cells.Item(6,1).Value = "test2"
$ws.Cells.Item(6,2).Value = "Test 2"
$ws.Cells.Item(7,1).Value = "options"
$ws.Cells.Item(7,2).Value = "OPTIONS"
$ws.Cells.Item(8,1).Value = "music"
$ws.Cells.Item(8,2).Value = "MUSIC"
$ws.Cells.Item(9,1).Value = "sound"
$ws.Cells.Item(9,2).Value = "SOUND"
$ws.Cells.Item(10,1).Value = "speech"
$ws.Cells.Item(10,2).Value = "SPEECH"
$ws.Cells.Item(11,1).Value = "on"
$ws.Cells.Item(11,2).Value = "ON"
$ws.Cells.Item(12,1).Value = "off"
$ws.Cells.Item(12,2).Value = "OFF"
$ws.Cells.Item(13,1).Value = "close"
$ws.Cells.Item(13,2).Value = "CLOSE"
$ws.Cells.Item(14,1).Value = "yes"
$ws.Cells.Item(14,2).Value = "YES"
$ws.Cells.Item(15,1).Value = "no"
$ws.Cells.Item(15,2).Value = "NO"
$ws.Cells.Item(16,1).Value = "testBodyCapsule"
$ws.Cells.Item(16,2).Value = "Capsule"
$ws.Cells.Item(17,1).Value = "testBodySphere"
$ws.Cells.Item(17,2).Value = "Sphere"
$ws.Cells.Item(18,1).Value = "categoryBody"
$ws.Cells.Item(18,2).Value = "Shape"
$ws.Cells.Item(19,1).Value = "categoryCellStructure"
$ws.Cells.Item(19,2).Value = "Structure"
$ws.Cells.Item(20,1).Value = "categoryMotility"
$ws.Cells.Item(20,2).Value = "Motility"
$ws.Cells.Item(21,1).Value = "categoryMetabolism"
$ws.Cells.Item(21,2).Value = "Metabolism"
$ws.Cells.Item(22,1).Value = "essentialNucleoid"
$ws.Cells.Item(22,2).Value = "Nucleoid"
$ws.Cells.Item(23,1).Value = "essentialRibosome"
$ws.Cells.Item(23,2).Value = "Ribosome"
$ws.Cells.Item(24,1).Value = "essentialPlasmid"
$ws.Cells.Item(24,2).Value = "Plasmid DNA"
$ws.Cells.Item(25,1).Value = "bodyBacillus"
$ws.Cells.Item(25,2).Value = "Bacillus"
$ws.Cells.Item(26,1).Value = "bodyCoccus"
$ws.Cells.Item(26,2).Value = "Coccus"
$ws.Cells.Item(27,1).Value = "bodyCoccobacillus"
$ws.Cells.Item(27,2).Value = "Coccobacillus"
$ws.Cells.Item(28,1).Value = "bodySpirillum"
$ws.Cells.Item(28,2).Value = "Spirillum"
$ws.Cells.Item(29,1).Value = "cellStructureThermophile"
$ws.Cells.Item(29,2).Value = "Thermophile"
$ws.Cells.Item(30,1).Value = "cellStructurePsychrophile"
$ws.Cells.Item(30,2).Value = "Psychrophile"
$ws.Cells.Item(31,1).Value = "cellStructureMethanogen"
$ws.Cells.Item(31,2).Value = "Methanogen"
$ws.Cells.Item(32,1).Value = "cellStructureHalophile"
$ws.Cells.Item(32,2).Value = "Halophile"
$ws.Cells.Item(33,1).Value = "motilityFlagellaMonotrichous"
$ws.Cells.Item(33,2).Value = "Monotrichous Flagella"
$ws.Cells.Item(34,1).Value = "motilityFlagellaLophotrichous"
$ws.Cells.Item(34,2).Value = "Lophotrichous Flagella"
$ws.Cells.Item(35,1).Value = "motilityFlagellaPeritrichous"
$ws.Cells.Item(35,2).Value = "Peritrichous Flagella"
$ws.Cells.Item(36,1).Value = "motilityFlagellaAmphitrichous"
$ws.Cells.Item(36,2).Value = "Amphitrichous Flagella"
$ws.Cells.Item(37,1).Value = "metabolismMethanotroph"
$ws.Cells.Item(37,2).Value = "Methanotroph"
$ws.Cells.Item(38,1).Value = "metabolismPhotoautotroph"
$ws.Cells.Item(38,2).Value = "Photoautotroph"
$ws.Cells.Item(39,1).Value = "metabolismOrganotroph"
$ws.Cells.Item(39,2).Value = "Organotroph"
$ws.Cells.Item(40,1).Value = "metabolismEndobiotic"
$ws.Cells.Item(40,2).Value = "Organotroph (Endobiotic)"

# Update the active selection to B16 (and drop the previous scrolled
# topLeftCell, matching the refreshed view state).
$null = $ws.Range("B16").Select()
